# Auto-generated edit script: update cryptos price/volume data
# (author commit: "Updated cryptos list on Wed Mar  8 18:34:55 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.021.68"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "1.553.90"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'290.85"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.3925"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").Value = "'0.3214"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'43.55"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").Value = "'0.07262"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "'1.074"
$ws.Range("E11").Value = "  -6.09%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.652"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'18.69"
$ws.Range("E14").Value = "  -8.22%  "
$ws.Range("D15").Value = "'0.00001126"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "'6.614"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "1.553.40"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'0.06589"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "'83.44"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'6.278"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "'15.50"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").Value = "'11.24"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D24").Value = "22.035.24"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "'2.375"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").Value = "'2.418"
$ws.Range("E26").Value = "  -6.26%  "
$ws.Range("D27").Value = "'148.69"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "'18.53"
$ws.Range("E28").Value = "  -4.23%  "
$ws.Range("D29").Value = "'4.890"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "1.727.80"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "'118.75"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("D32").Value = "'0.9961"
$ws.Range("E32").Value = "  -8.07%  "
$ws.Range("D33").Value = "'5.765"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").Value = "'0.08304"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "'1.621"
$ws.Range("E35").Value = "  -15.70%  "
$ws.Range("D36").Value = "'8.973"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("D37").Value = "'0.02256"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").Value = "'0.06077"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("D39").Value = "'5.091"
$ws.Range("E39").Value = "  -5.70%  "
$ws.Range("D40").Value = "'1.208"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Value = "'0.2039"
$ws.Range("E41").Value = "  -6.00%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "'10.66"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("D44").Value = "'0.5808"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.751"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.93"
$ws.Range("E46").Value = "  -6.51%  "
$ws.Range("D47").Value = "'0.5570"
$ws.Range("E47").Value = "  -5.94%  "
$ws.Range("D48").Value = "'118.25"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").Value = "'1.894"
$ws.Range("E49").Value = "  -4.63%  "
$ws.Range("D50").Value = "'1.132"
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("D51").Value = "'0.06820"
$ws.Range("E51").Value = "  -3.87%  "
